$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A105:A143").Value = 145.164347804412
